# Update cryptocurrency price/volume data in the worksheet
# reflecting the latest scrape values (GitHub Actions scheduled update).
#
# Note: several "Price" values in column D are strings that look like
# plain numbers (e.g. "524.55"); Excel auto-converts such text into a
# floating-point number on assignment, which introduces binary rounding
# noise (e.g. 524.54999999999995) and loses the exact original text.
# To preserve the exact literal text, we force a Text number format
# before writing the value and then restore the default "Normal" style
# so no lasting formatting change is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.455.77"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "3.106.92"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.00%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.101.84"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  +2.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "3.638.59"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "57.545.73"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "3.110.16"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "347.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("E26").Value = "  -1.76%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "0.0₃0907"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("E33").Value = "  -7.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.37%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.17%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.70%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.65%  "
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "3.142.15"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.346.12"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +3.34%  "
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.70%  "
